$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106, shifting existing rows 106-108 down to 107-109
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with the new data record
$ws.Range("A106").Value = 9
$ws.Range("B106").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C106").Value = 'Metropolitana'
$ws.Range("D106").Value = 44628
$ws.Range("E106").Value = 13
$ws.Range("F106").Value = 100112022
$ws.Range("G106").Value = 'Arveja Verde'
$ws.Range("H106").Value = 'Sin especificar'
$ws.Range("I106").Value = 'Primera'
$ws.Range("J106").Value = 43
$ws.Range("K106").Value = 24000
$ws.Range("L106").Value = 26000
$ws.Range("M106").Value = 25023
$ws.Range("N106").Value = '$/malla 25 kilos'
$ws.Range("O106").Value = 'Carahue'
$ws.Range("P106").Value = 1001
$ws.Range("Q106").Value = 25
$ws.Range("R106").Value = 'Hortaliza'
